# Insert a new price record at row 65 (weekly Fruta/Hortalizas update),
# pushing the existing rows 65:170 down to 66:171.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("65:65").Insert()

# Populate the new row 65 with the latest observation.
$ws.Cells.Item(65, 1).Value = 5
$ws.Cells.Item(65, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(65, 3).Value = "Maule"
$ws.Cells.Item(65, 4).Value = 45036
$ws.Cells.Item(65, 5).Value = 7
$ws.Cells.Item(65, 6).Value = "Fruta"
$ws.Cells.Item(65, 7).Value = 100108
$ws.Cells.Item(65, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(65, 9).Value = 100108002
$ws.Cells.Item(65, 10).Value = "Mango"
$ws.Cells.Item(65, 11).Value = "Sin especificar"
$ws.Cells.Item(65, 12).Value = "Primera"
$ws.Cells.Item(65, 13).Value = 248
$ws.Cells.Item(65, 14).Value = 7000
$ws.Cells.Item(65, 15).Value = 7000
$ws.Cells.Item(65, 16).Value = 7000
$ws.Cells.Item(65, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(65, 18).Value = "Perú"
$ws.Cells.Item(65, 19).Value = 1750
$ws.Cells.Item(65, 20).Value = 4
